$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FetHead")

# Move the existing total row (row 22 -> row 23) down to make room for the
# new "DIP Switch" line item.
$ws.Rows.Item(22).Insert()

# New row 21: DIP Switch / Address selector line item
$ws.Range("B21").Value = "DIP Switch"
$ws.Range("D21").Value = 1
$ws.Range("G21").Value = "Address selector"

# Update selection to match the saved workbook state
$ws.Range("H19").Select()
